# cierre 5 May 22
# ------------------------------------------------------------------
# This script reproduces, via the Excel COM object model, the edits
# captured in the target diff:
#   1. Sheet "REMISIONES   ABRIL  2 0 2 2    " (credit ledger for April):
#      six blank trailing rows (old rows 60:65) are deleted and the
#      folio number left in B59 is cleared.
#   2. Sheet "Hoja4" is filled in with the "cierre" (closing) table of
#      bank-deposit entries that used to be an empty sheet, and becomes
#      the active sheet/tab.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ============================================================
# 1) REMISIONES ABRIL 2 0 2 2 -- delete the 6 now-unused rows
# ============================================================
$wsAbril = $wb.Worksheets.Item("REMISIONES   ABRIL  2 0 2 2    ")
$wsAbril.Activate()

# Delete rows 60:65 (whole rows) -- everything below shifts up by 6,
# shared formulas / SUM ranges / merged cells / drawings anchored to
# those rows move automatically.
$wsAbril.Range("A60:A65").EntireRow.Delete() | Out-Null

# Row 59 keeps its place but its folio number (340) is removed.
$wsAbril.Range("B59").ClearContents() | Out-Null

Write-Host "Abril sheet: rows 60:65 deleted, B59 cleared"

# ============================================================
# 2) Hoja4 -- "cierre 5 May 22" bank deposits closing table
# ============================================================
$ws = $wb.Worksheets.Item("Hoja4")

# --- column widths -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.85546875
$ws.Columns.Item(2).ColumnWidth = 16.140625
$ws.Columns.Item(3).ColumnWidth = 17.140625
$ws.Columns.Item(4).ColumnWidth = 18.140625
$ws.Columns.Item(5).ColumnWidth = 17.42578125

# --- title -----------------------------------------------------------
$ws.Range("B1").Value = "ABASTOS DE 4 CANRES  HERRADURA "
$ws.Rows.Item(1).RowHeight = 21

Write-Host "Hoja4 started"

